$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "89.495.58"
$ws.Range("E2").Value = "  -1.73%  "

# Row 3
$ws.Range("D3").Value = "3.100.10"
$ws.Range("E3").Value = "  -2.27%  "

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "'213.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.78%  "

# Row 6
$ws.Range("D6").Value = "'622.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.94%  "

# Row 7
$ws.Range("E7").Value = "  -4.91%  "

# Row 8
$ws.Range("E8").Value = "  +16.02%  "

# Row 9
$ws.Range("E9").Value = "  +0.13%  "

# Row 10
$ws.Range("D10").Value = "3.098.97"
$ws.Range("E10").Value = "  -2.16%  "

# Row 11
$ws.Range("E11").Value = "  +8.46%  "

# Row 12
$ws.Range("E12").Value = "  +1.52%  "

# Row 13
$ws.Range("D13").Value = "'0.0000242"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.78%  "

# Row 14
$ws.Range("E14").Value = "  +0.59%  "

# Row 15
$ws.Range("D15").Value = "89.195.73"
$ws.Range("E15").Value = "  -1.54%  "

# Row 16
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.678.65"
$ws.Range("E16").Value = "  -1.89%  "

# Row 17
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'32.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.14%  "

# Row 18
$ws.Range("D18").Value = "3.099.92"
$ws.Range("E18").Value = "  -2.21%  "

# Row 19
$ws.Range("D19").Value = "'3.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.23%  "

# Row 20
$ws.Range("D20").Value = "'0.0000215"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.94%  "

# Row 21
$ws.Range("D21").Value = "'13.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.22%  "

# Row 22
$ws.Range("D22").Value = "'426.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.01%  "

# Row 23
$ws.Range("D23").Value = "'8.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.34%  "

# Row 24
$ws.Range("E24").Value = "  -0.29%  "

# Row 25
$ws.Range("D25").Value = "'5.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.84%  "

# Row 26
$ws.Range("D26").Value = "'12.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.73%  "

# Row 27
$ws.Range("D27").Value = "'83.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.01%  "

# Row 28
$ws.Range("D28").Value = "3.257.65"
$ws.Range("E28").Value = "  -2.33%  "

# Row 29
$ws.Range("E29").Value = "  +0.10%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.32%  "

# Row 31
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "'0.164"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.81%  "

# Row 32
$ws.Range("E32").Value = "  -1.39%  "

# Row 33
$ws.Range("D33").Value = "'511.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.17%  "

# Row 34
$ws.Range("E34").Value = "  -6.87%  "

# Row 35
$ws.Range("D35").Value = "'6.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.48%  "

# Row 36
$ws.Range("D36").Value = "'1.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.26%  "

# Row 37
$ws.Range("E37").Value = "  -3.61%  "

# Row 38
$ws.Range("D38").Value = "'22.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "

# Row 39
$ws.Range("D39").Value = "'22.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.42%  "

# Row 40
$ws.Range("D40").Value = "'0.129"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.30%  "

# Row 41
$ws.Range("E41").Value = "  +0.29%  "

# Row 42
$ws.Range("E42").Value = "  +0.02%  "

# Row 43
$ws.Range("E43").Value = "  -0.61%  "

# Row 44
$ws.Range("E44").Value = "  -3.62%  "

# Row 45
$ws.Range("D45").Value = "'0.134"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.96%  "

# Row 46
$ws.Range("D46").Value = "'145.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "

# Row 47
$ws.Range("D47").Value = "'0.0710"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.35%  "

# Row 48
$ws.Range("D48").Value = "'43.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.64%  "

# Row 49
$ws.Range("E49").Value = "  +2.28%  "

# Row 50
$ws.Range("D50").Value = "'160.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.60%  "

# Row 51
$ws.Range("D51").Value = "'0.707"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.05%  "
